# Steel_Pan_Head_Screws_with_Internal-Tooth_Lock_Washer.xlsx
#
# The sheet gains two new rows at the very top:
#   - new row 1: a header of column indices (0..13) using the old header's
#     bold/bordered style
#   - new row 2: a mostly blank row, except column E which reads "Washer"
#   - the old header row ("Lg.", "Threading", ... ) is pushed down to row 3,
#     loses its bold/border styling, and its M/N cells (which used to hold
#     "thread_size" / "material_surface") are cleared
#   - every former data row shifts down by two rows (old row 2 -> row 4,
#     ..., old row 66 -> row 68)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 1; this pushes everything (including the
# header row) down by two rows.
$ws.Range("A1:A2").EntireRow.Insert()

# The old header row (with its bold/bordered style) now lives at row 3.
# Copy its formatting up to the new row 1 before we strip it from row 3.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)  # xlPasteFormats

# Populate new row 1 with the numeric column index header (0-based).
for ($i = 0; $i -lt 14; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $i
}

# New row 2 is blank except for "Washer" in column E.
$ws.Range("E2").Value = "Washer"

# Strip the bold/bordered styling that row 3 inherited when it shifted down
# from row 1 (the target layout has the plain, unstyled header text here).
$ws.Range("A3:N3").ClearFormats()

# Row 3's M/N cells used to carry "thread_size" / "material_surface" - those
# are no longer present in the new layout.
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
